$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet "Correlation matrix": the config now carries only 2 HW inputs and
# 2 HD inputs (one config file shared across gadgets) instead of 8 each, so
# the matrix shrinks from 16 data columns down to 4, and every value is
# recomputed from the smaller run.
# ---------------------------------------------------------------------------
$wsCorr = $wb.Worksheets.Item("Correlation matrix")

# Drop the now-unused HW input2..7 and HD input2..7 columns (D..I and L..Q),
# shifting HD input0/HD input1 (currently J/K) left into D/E.
$wsCorr.Range("L1:Q4").EntireColumn.Delete()
$wsCorr.Range("D1:I4").EntireColumn.Delete()

# Relabel / rewrite the surviving header + data cells with the new run's
# values.
$wsCorr.Range("B1").Value = "HW input0"
$wsCorr.Range("C1").Value = "HW input1"
$wsCorr.Range("D1").Value = "HD input0"
$wsCorr.Range("E1").Value = "HD input1"

$wsCorr.Range("A2").Value = "no delays"
$wsCorr.Range("B2").Value = 0.003937007874015723
$wsCorr.Range("C2").Value = 0.003937007874015723
$wsCorr.Range("D2").Value = -0.003937007874015723
$wsCorr.Range("E2").Value = -0.003937007874015723

$wsCorr.Range("A3").Value = "gate delays"
$wsCorr.Range("B3").Value = 0.1854742807555382
$wsCorr.Range("C3").Value = 0.1854742807555382
$wsCorr.Range("D3").Value = 0.1755634413258522
$wsCorr.Range("E3").Value = 0.1755634413258529

$wsCorr.Range("A4").Value = "gate+inputs delay"
$wsCorr.Range("B4").Value = 0.1854742807555382
$wsCorr.Range("C4").Value = 0.1854742807555382
$wsCorr.Range("D4").Value = 0.1755634413258522
$wsCorr.Range("E4").Value = 0.1755634413258529

# ---------------------------------------------------------------------------
# Sheet "Toggles no del": one fewer data point (2 instead of 3).
# ---------------------------------------------------------------------------
$wsNoDel = $wb.Worksheets.Item("Toggles no del")
$wsNoDel.Range("A4").EntireRow.Delete()
$wsNoDel.Range("B2").Value = 127
$wsNoDel.Range("B3").Value = 128

$chart1 = $wsNoDel.ChartObjects().Item(1).Chart
$chart1.SeriesCollection().Item(1).Formula = "=SERIES(,,'Toggles no del'!`$B`$2:`$B`$3,1)"

# ---------------------------------------------------------------------------
# Sheet "Toggles del": one fewer data point (3 instead of 4).
# ---------------------------------------------------------------------------
$wsDel = $wb.Worksheets.Item("Toggles del")
$wsDel.Range("A5").EntireRow.Delete()
$wsDel.Range("B2").Value = 79
$wsDel.Range("B3").Value = 128
$wsDel.Range("B4").Value = 48

$chart2 = $wsDel.ChartObjects().Item(1).Chart
$chart2.SeriesCollection().Item(1).Formula = "=SERIES(,,'Toggles del'!`$B`$2:`$B`$4,1)"

# ---------------------------------------------------------------------------
# Sheet "Toggles input del": same reduction, same new values as "Toggles del".
# ---------------------------------------------------------------------------
$wsInputDel = $wb.Worksheets.Item("Toggles input del")
$wsInputDel.Range("A5").EntireRow.Delete()
$wsInputDel.Range("B2").Value = 79
$wsInputDel.Range("B3").Value = 128
$wsInputDel.Range("B4").Value = 48

$chart3 = $wsInputDel.ChartObjects().Item(1).Chart
$chart3.SeriesCollection().Item(1).Formula = "=SERIES(,,'Toggles input del'!`$B`$2:`$B`$4,1)"
